$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.756.84"
$ws.Range("E2").Value = "  +1.51%  "
$ws.Range("D3").Value = "1.648.51"
$ws.Range("E3").Value = "  -0.32%  "
$ws.Range("E4").Value = "  +0.07%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "213.54"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.20%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.534"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +3.35%  "
$ws.Range("E7").Value = "  +0.12%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "23.22"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -0.91%  "
$ws.Range("E9").Value = "  +0.18%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.0615"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +0.18%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0890"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.81%  "
$ws.Range("D12").Value = "1.881.26"
$ws.Range("E12").Value = "  -0.20%  "
$ws.Range("D13").Value = "1.646.57"
$ws.Range("E13").Value = "  -0.30%  "
$ws.Range("E14").Value = "  -0.73%  "
$ws.Range("E15").Value = "  -0.81%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "64.39"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -1.70%  "
$ws.Range("D17").Value = "27.726.70"
$ws.Range("E17").Value = "  +1.47%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "231.82"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +0.79%  "
$ws.Range("E19").Value = "  -0.13%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "7.67"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +3.57%  "
$ws.Range("E22").Value = "  -0.90%  "
$ws.Range("E23").Value = "  +7.71%  "
$ws.Range("E24").Value = "  -3.75%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "149.96"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +1.80%  "
$ws.Range("E26").Value = "  -1.24%  "
$ws.Range("E27").Value = "  +0.62%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "15.69"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -0.78%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("E30").Value = "  +0.30%  "
$ws.Range("E31").Value = "  -1.71%  "
$ws.Range("E32").Value = "  +0.27%  "
$ws.Range("D33").Value = "1.444.69"
$ws.Range("E33").Value = "  +1.18%  "
$ws.Range("E34").Value = "  +1.26%  "
$ws.Range("E35").Value = "  +1.35%  "
$ws.Range("E37").Value = "  +0.45%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.888"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -2.08%  "
$ws.Range("E39").Value = "  -0.43%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.886"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +12.48%  "
$ws.Range("E41").Value = "  -1.94%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "5.71"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +2.52%  "
$ws.Range("E43").Value = "  +0.01%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "66.52"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +2.36%  "
$ws.Range("E45").Value = "  -0.30%  "
$ws.Range("E46").Value = "  +2.06%  "
$ws.Range("D47").Value = "1.790.57"
$ws.Range("E47").Value = "  -0.14%  "
$ws.Range("E48").Value = "  +3.56%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "86.50"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -1.73%  "
$ws.Range("E50").Value = "  +1.76%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.0996"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -1.66%  "
